# Daily attendance processing - 2026-01-02 06:43:59
# Re-order the "Recorded By" (column G) comma-separated list so that the
# capitalized "System" token is moved to immediately follow a lowercase
# "system" token if one is present, otherwise to the front of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notmatch "System") { continue }

    $tokens = $val -split ", "

    # remove the exact "System" token(s) -- NOTE: this runtime's -ceq/-cne/
    # -cmatch/-clike operators are NOT actually case-sensitive, so use
    # [string].Equals(), which performs an ordinal (case-sensitive) compare.
    $withoutSystem = @()
    $foundSystem = $false
    foreach ($tok in $tokens) {
        if ($tok.Equals("System")) {
            $foundSystem = $true
        } else {
            $withoutSystem += $tok
        }
    }

    if (-not $foundSystem) {
        # no exact "System" token found (case-sensitive) - leave untouched
        continue
    }

    $lowerIdx = -1
    for ($i = 0; $i -lt $withoutSystem.Count; $i++) {
        if ($withoutSystem[$i].Equals("system")) {
            $lowerIdx = $i
            break
        }
    }

    if ($lowerIdx -ge 0) {
        $newTokens = @()
        for ($i = 0; $i -le $lowerIdx; $i++) { $newTokens += $withoutSystem[$i] }
        $newTokens += "System"
        for ($i = $lowerIdx + 1; $i -lt $withoutSystem.Count; $i++) { $newTokens += $withoutSystem[$i] }
    } else {
        $newTokens = @("System") + $withoutSystem
    }

    $newVal = $newTokens -join ", "
    $cell.Value2 = $newVal
}
